# Slide 2 / shape 1 is the empty "1 Título" (title) placeholder. The author
# typed the question that titles this slide into it, left-aligned the
# paragraph and PowerPoint auto-shrank the text to fit the placeholder.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$title = $s.Shapes.Item(1)
$tf = $title.TextFrame
$tr = $tf.TextRange

# Type the text in (InsertAfter on the still-empty range keeps the
# paragraph's existing trailing end-of-paragraph run properties instead of
# fabricating a brand new, formatting-less paragraph).
[void]$tr.InsertAfter("Que es la Ofimática?")

# Mark the new run as Spanish, matching the rest of the deck.
$tf.TextRange.LanguageID = "es-ES"

# Left align the paragraph (title placeholder default is centered).
$tf.TextRange.ParagraphFormat.Alignment = 1      # ppAlignLeft -> <a:pPr algn="l"/>

# Let PowerPoint shrink the text to fit the title placeholder.
$tf.AutoSize = 2                                 # ppAutoSizeTextToFitShape -> <a:normAutofit/>
